$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7841
$ws1.Range("F4").Value = 223
$ws1.Range("F5").Value = 58
$ws1.Range("F6").Value = 574
$ws1.Range("F7").Value = 1190
$ws1.Range("F10").Value = 177

# Sheet "全部类型" (all types list) - mirrors the same rows, shifted by one
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7841
$ws4.Range("F4").Value = 223
$ws4.Range("F5").Value = 58
$ws4.Range("F6").Value = 574
$ws4.Range("F7").Value = 1190
$ws4.Range("F11").Value = 177
